$wb = $excel.ActiveWorkbook

# --- Prerequisites sheet ---
$wsPre = $wb.Worksheets.Item("Prerequisites")
$wsPre.Cells.Item(2, 3).Value = "Prerequisites"
$wsPre.Cells.Item(3, 3).Value = "That the web page is working.
That the manual tests have been successful.
Intellij IDEA Community
Gradle v6.8
JDK v11
Chrome
Chromedriver
Plugins:
  Cucumber for Java
  Gherkin
Have internet"

$wsPre.Columns.Item(3).ColumnWidth = 37.7109375

$headerPre = $wsPre.Range("C2")
$headerPre.Font.Bold = $true
$headerPre.Font.Size = 12
$headerPre.Interior.Pattern = -4124
$headerPre.Interior.PatternColorIndex = -4105
$headerPre.Interior.ThemeColor = 9
$headerPre.Interior.TintAndShade = 0.39997558519241921
$headerPre.HorizontalAlignment = -4108
$headerPre.Borders.LineStyle = 1
$headerPre.Borders.Weight = 2
$wsPre.Rows.Item(2).RowHeight = 15.75

$bodyPre = $wsPre.Range("C3")
$bodyPre.HorizontalAlignment = -4131
$bodyPre.VerticalAlignment = -4108
$bodyPre.WrapText = $true
$bodyPre.Borders.LineStyle = 1
$bodyPre.Borders.Weight = 2
$wsPre.Rows.Item(3).RowHeight = 180

$wsPre.Range("A1").Select()
$wsPre.Application.ActiveWindow.Zoom = 130
$wsPre.Application.ActiveWindow.DisplayGridlines = $false
$wsPre.Range("C3").Select()

# --- Limitations sheet ---
$wsLim = $wb.Worksheets.Item("Limitations")
$wsLim.Cells.Item(2, 3).Value = "Limitations"
$wsLim.Cells.Item(3, 3).Value = "We have 30mbps navigation."

$wsLim.Columns.Item(3).ColumnWidth = 38.140625

$headerLim = $wsLim.Range("C2")
$headerLim.Font.Bold = $true
$headerLim.Font.Size = 12
$headerLim.Interior.Pattern = -4124
$headerLim.Interior.PatternColorIndex = -4105
$headerLim.Interior.ThemeColor = 9
$headerLim.Interior.TintAndShade = 0.39997558519241921
$headerLim.HorizontalAlignment = -4108
$headerLim.Borders.LineStyle = 1
$headerLim.Borders.Weight = 2
$wsLim.Rows.Item(2).RowHeight = 15.75

$bodyLim = $wsLim.Range("C3")
$bodyLim.HorizontalAlignment = -4131
$bodyLim.VerticalAlignment = -4160
$bodyLim.WrapText = $true
$bodyLim.Borders.LineStyle = 1
$bodyLim.Borders.Weight = 2
$wsLim.Rows.Item(3).RowHeight = 24.75

$wsLim.Range("A1").Select()
$wsLim.Application.ActiveWindow.Zoom = 130
$wsLim.Application.ActiveWindow.DisplayGridlines = $false
$wsLim.Range("C2").Select()

# --- Make Assumptions the active tab ---
$wsAssumptions = $wb.Worksheets.Item("Assumptions")
$wsAssumptions.Activate()
